# "Generate Report for Handoff"
#
# The two source files that used to be independently handed-back
# (3ca43955-...md and bbeb3347-...md) are replaced by a freshly
# generated handoff report for two new source files
# (09adf306-...md and ffffec10c4c6-...md), both bundled into a single
# combined handoff package (...2bc493a2...zh-cn.xlf / ...de-de.xlf).
# Status flips from "Handed back: in sync with en-US" to
# "Ready for handoff", the handoff datetime advances, the handback
# datetime resets to the zero-date (nothing handed back yet), and the
# now-unused "Latest Target File" / "Latest Handback File" columns are
# cleared out for the two data rows.

$wb = $excel.ActiveWorkbook

$oldFile1 = "3ca43955-bd38-4955-aafa-69ccbe78ff67.md"
$oldFile2 = "bbeb3347-5d2c-450a-b41a-f9167885289b.md"

$newFile1 = "09adf306-b146-44f5-90ff-b90e86020fba.md"
$newFile2 = "ffffec10c4c6-7925-49ef-9628-fd17d695b3e2.md"

$statusReady = "Ready for handoff"

$zhHandoffXlf = "09adf306-b146-44f5-90ff-b90e86020fba.2bc493a2c28fc25a8ed57886914c2e8491c2ed5c.zh-cn.xlf"
$deHandoffXlf = "09adf306-b146-44f5-90ff-b90e86020fba.2bc493a2c28fc25a8ed57886914c2e8491c2ed5c.de-de.xlf"

$zhHandoffDatetime = "2016-03-08 10:43:55"
$deHandoffDatetime = "2016-03-08 10:43:59"
$zeroDatetime = "0001-01-01 00:00:00"

$srcRepoCommit = "04e13ba9491c405a4dd4a5fd6421fdfe320c02cc"

# ----------------------------------------------------------------
# Overview sheet
# ----------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Hyperlinks.Delete()

$ov.Range("B2").Value = $statusReady
$ov.Range("C2").Value = $statusReady
$ov.Range("B3").Value = $statusReady
$ov.Range("C3").Value = $statusReady

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcRepoCommit/e2e/$newFile1", "", "", $newFile1)
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcRepoCommit/e2e/$newFile2", "", "", $newFile2)
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcRepoCommit/.localization-config", "", "", ".localization-config")

$ov.Range("A2").Style = 1
$ov.Range("A3").Style = 1
$ov.Range("A4").Style = 1

# ----------------------------------------------------------------
# Per-locale sheets (zh-cn, de-de)
# ----------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; HandoffXlf = $zhHandoffXlf; HandoffDatetime = $zhHandoffDatetime; HandoffRepo = "olhandoff"; HandoffCommit = "0dcfe663ff287638f1e7daf8903e2c9222e7438a" },
    @{ Sheet = "de-de"; HandoffXlf = $deHandoffXlf; HandoffDatetime = $deHandoffDatetime; HandoffRepo = "olhandoff"; HandoffCommit = "39c69f6d5f07303fc309ba23d00ae24b55eacf66" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)
    $locale = $loc.Sheet

    # Clear every hyperlink on the sheet so stale rIds for the removed
    # E/F columns and renamed files don't linger.
    $ws.Hyperlinks.Delete()

    # Row 2: first source file, now pointing at the shared handoff package
    $ws.Range("B2").Value = $statusReady
    $ws.Range("C2").Value = $loc.HandoffXlf
    $ws.Range("D2").Value = $loc.HandoffDatetime
    $ws.Range("E2").Clear()
    $ws.Range("F2").Clear()
    $ws.Range("G2").Value = $zeroDatetime
    $ws.Range("H2").Value = "Include"

    # Row 3: second source file, same shared handoff package
    $ws.Range("B3").Value = $statusReady
    $ws.Range("C3").Value = $loc.HandoffXlf
    $ws.Range("D3").Value = $loc.HandoffDatetime
    $ws.Range("E3").Clear()
    $ws.Range("F3").Clear()
    $ws.Range("G3").Value = $zeroDatetime
    $ws.Range("H3").Value = "Include"

    # Row 4 (.localization-config) is untouched content-wise.

    $handoffXlfUrl = "https://github.com/OpenLocalizationTestOrg/$($loc.HandoffRepo)/blob/$($loc.HandoffCommit)/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/$($loc.HandoffXlf)"

    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcRepoCommit/e2e/$newFile1", "", "", $newFile1)
    $ws.Hyperlinks.Add($ws.Range("C2"), $handoffXlfUrl, "", "", $loc.HandoffXlf)
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcRepoCommit/e2e/$newFile2", "", "", $newFile2)
    $ws.Hyperlinks.Add($ws.Range("C3"), $handoffXlfUrl, "", "", $loc.HandoffXlf)
    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcRepoCommit/.localization-config", "", "", ".localization-config")

    $ws.Range("A2").Style = 1
    $ws.Range("C2").Style = 1
    $ws.Range("A3").Style = 1
    $ws.Range("C3").Style = 1
    $ws.Range("A4").Style = 1
}
